# Edit summary
# ------------
# 1) Slide 6 contains a 4-column table. Its table style (design-gallery
#    style) is switched from the deck's custom "Table_0" style to the
#    PowerPoint built-in style {2D3AA692-0FB0-4110-A5C7-DEBA2E7F1FDE}.
#    Table styles can't be assigned through the Style property directly
#    (PowerPoint raises "Table styles cannot be assigned through a
#    property"), so we use Table.ApplyStyle with the brace-GUID id.
#
# 2) The presentation's theme was switched away from the custom
#    "Integral" design back to the default Office theme palette (the
#    deck keeps using the Integral fonts/effects scheme, only the
#    colour scheme changes). We push the 12 standard Office theme
#    colours through the live ThemeColorScheme so the underlying theme
#    part picks up the default Office palette.

$p = $ppt.ActivePresentation

# --- 1) Table style -------------------------------------------------
$slide = $p.Slides.Item(6)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{2D3AA692-0FB0-4110-A5C7-DEBA2E7F1FDE}")

# --- 2) Theme colours -------------------------------------------------
# Index order (MsoThemeColorSchemeIndex): 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1-6, 11 hlink, 12 folHlink. Values below are the default
# Office theme colours (RGB packed as r + g*256 + b*65536, matching the
# VBA RGB() helper).
$officeColors = @{
    1  = 0x000000   # dk1      000000
    2  = 0xFFFFFF   # lt1      FFFFFF
    3  = 0x6A5444   # dk2      44546A
    4  = 0xE6E6E7   # lt2      E7E6E6
    5  = 0xD59B5B   # accent1  5B9BD5
    6  = 0x317DED   # accent2  ED7D31
    7  = 0xA5A5A5   # accent3  A5A5A5
    8  = 0x00C0FF   # accent4  FFC000
    9  = 0xC47244   # accent5  4472C4
    10 = 0x47AD70   # accent6  70AD47
    11 = 0xC16305   # hlink    0563C1
    12 = 0x724F95   # folHlink 954F72
}

$cs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $cs.Colors($i).RGB = $officeColors[$i]
}
